# Refresh the cryptocurrency snapshot (prices / 1h volume change, plus the
# Arweave / Bittensor row swap) to match the latest GitHub Actions run.
#
# Cells hold plain text (prices use "." as a thousands separator, e.g.
# "70.190.27", and percentages are padded strings like "  -0.19%  "), so
# every cell is forced to text format ("@") before its value is written --
# otherwise Excel auto-coerces values such as "1.00" or "16.40" into the
# numbers 1 / 16.4, silently dropping the trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.190.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.747.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.31%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.65%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.745.68"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.35%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.62%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.48%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.48%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.482"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.47%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.15"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.66%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000255"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.365.60"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.743.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.45%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.238.82"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.42%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.58"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.58%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "505.43"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.40"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.86%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.58%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.56"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.41%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.84%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.02"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.14"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.61%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.30%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.47"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.93"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.34%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.89"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.52%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.56"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.38%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.93%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.14%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.10"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.88%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.348"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.55%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.54%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.21"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +15.53%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.08"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.44%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.87"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.63%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "429.23"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.41%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.58"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.75%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.25%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.969.81"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0362"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.92%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.26"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.38%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.41"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.40%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.87%  "
